# DEA_input_dataset.xlsx - "modify dea input sample, finish all features"
#
# Appends 25 new data rows (96-120) to the DEA_input_data sheet, re-using
# the same A:E layout as the existing rows (Bus Line/Customer id, Hours,
# Miles, Bus Count, Customer). String-valued ids go through the shared
# string table automatically when assigned via .Value; purely numeric ids
# (701, 703, 628, 704, 750, 720) are written as plain numbers, matching
# the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 96;  A = "2X";   B = 450.14666666699998; C = 772.99052977199995; D = 12; E = 19.091281285000001 },
    @{ Row = 97;  A = "F590"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 98;  A = "F518"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 99;  A = "F514"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 100; A = "F618"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 101; A = "F504"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 102; A = "F534"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 103; A = "F638"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 104; A = "F522"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 105; A = "F400"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 106; A = "F402"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 107; A = "F94";  B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 108; A = "F556"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 109; A = "F570"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 110; A = "F547"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 111; A = "F546"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 112; A = "F401"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 113; A = "F578"; B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 114; A = 701;    B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 115; A = 703;    B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 116; A = 628;    B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 117; A = 704;    B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 118; A = 750;    B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 119; A = 720;    B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 },
    @{ Row = 120; A = "35M";  B = 269.36;              C = 400.12379214399999; D = 4;  E = 101.52455476199999 }
)

foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
}

# Match the author's final on-screen view: scrolled down with C126 selected
# (scroll position itself isn't persisted by this host, but the selection is).
try {
    $excel.ActiveWindow.ScrollRow = 107
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("C126").Select()
